$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.668.68'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '2.487.28'
$ws.Range('E3').Value = '  -1.56%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '531.41'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.08'
$ws.Range('E6').Value = '  -2.57%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '2.506.79'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').Value = '2.931.93'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '23.11'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').Value = '58.628.58'
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000139'
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('D18').Value = '2.502.32'
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.97'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('E20').Value = '  +0.36%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '323.51'
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.84'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '64.11'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('E26').Value = '  -0.62%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('E28').Value = '  -3.59%  '
$ws.Range('D29').Value = '0.0₃0768'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.66'
$ws.Range('E30').Value = '  -0.91%  '
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '166.74'
$ws.Range('E32').Value = '  +3.60%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.14'
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.37'
$ws.Range('E35').Value = '  -4.70%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '18.45'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('E37').Value = '  -3.59%  '
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.811'
$ws.Range('E40').Value = '  +1.07%  '
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.24'
$ws.Range('E42').Value = '  -1.52%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '277.52'
$ws.Range('E43').Value = '  -3.14%  '
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.599'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.86'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '127.91'
$ws.Range('E47').Value = '  +3.30%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0923'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0511'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '17.26'
$ws.Range('E51').Value = '  -1.44%  '
